$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Alief Faza Rizqi Adi Jaya): mark Tugas 2 (D) and Tugas 4 (F) as completed
$ws.Range("D4").Value = "ü"
$ws.Range("D4").Font.Name = "Wingdings"
$ws.Range("D4").Font.Size = 12

$ws.Range("F4").Value = "ü"
$ws.Range("F4").Font.Name = "Wingdings"
$ws.Range("F4").Font.Size = 12

# Row 19 (Muhammad Farhan): mark Tugas 2 (D) and Tugas 4 (F) as completed
$ws.Range("D19").Value = "ü"
$ws.Range("D19").Font.Name = "Wingdings"
$ws.Range("D19").Font.Size = 12

$ws.Range("F19").Value = "ü"
$ws.Range("F19").Font.Name = "Wingdings"
$ws.Range("F19").Font.Size = 12
